$wb = $excel.ActiveWorkbook

# Sheet1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1473   # was 1470
$ws.Range("F4").Value = 799   # was 796
$ws.Range("F8").Value = 7166   # was 7163
$ws.Range("F10").Value = 5217   # was 5214
$ws.Range("F15").Value = 8471   # was 8461
$ws.Range("F17").Value = 1122   # was 1121
$ws.Range("F18").Value = 827   # was 825
$ws.Range("F27").Value = 1582   # was 1580
$ws.Range("F29").Value = 831   # was 829
$ws.Range("F30").Value = 1816   # was 1814
$ws.Range("F32").Value = 2179   # was 2178
$ws.Range("F35").Value = 1390   # was 1389
$ws.Range("F36").Value = 65   # was 64
$ws.Range("F38").Value = 766   # was 765
$ws.Range("F40").Value = 2876   # was 2875
$ws.Range("F41").Value = 3946   # was 3945
$ws.Range("F42").Value = 181   # was 180
$ws.Range("F48").Value = 144   # was 143

# Sheet2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F17").Value = 0   # was 40

# Sheet3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 4952   # was 4948

# Sheet4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 4952   # was 4948
$ws.Range("F6").Value = 1473   # was 1470
$ws.Range("F7").Value = 799   # was 796
$ws.Range("F13").Value = 5217   # was 5214
$ws.Range("F19").Value = 1122   # was 1121
$ws.Range("F20").Value = 827   # was 825
$ws.Range("F29").Value = 1582   # was 1580
$ws.Range("F31").Value = 831   # was 829
$ws.Range("F32").Value = 1816   # was 1814
$ws.Range("F34").Value = 2179   # was 2178
$ws.Range("F39").Value = 766   # was 765
$ws.Range("F43").Value = 3946   # was 3945
$ws.Range("F44").Value = 181   # was 180
$ws.Range("F48").Value = 144   # was 143
